$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1374435.8
$ws.Range("I70").Value = 5166.6665
$ws.Range("J70").Value = 2121309.8
$ws.Range("K70").Value = 15499.9995
$ws.Range("L70").Value = 6363929.399999999
$ws.Range("M70").Value = -15229.9995
$ws.Range("N70").Value = -6364469.399999999
$ws.Range("H73").Value = 1374435.8
$ws.Range("I73").Value = 5166.6665
$ws.Range("J73").Value = 2121309.8
$ws.Range("K73").Value = 15499.9995
$ws.Range("L73").Value = 6363929.399999999
$ws.Range("M73").Value = -14563.9995
$ws.Range("N73").Value = -6365801.399999999
$ws.Range("H98").Value = 12090.786
$ws.Range("I98").Value = 14597.637
$ws.Range("K98").Value = 14597.637
$ws.Range("M98").Value = -13099.637
$ws.Range("H106").Value = 3272.7273
$ws.Range("J106").Value = 3375
$ws.Range("L106").Value = 3375
$ws.Range("N106").Value = -4637
$ws.Range("H107").Value = 1330.4117
$ws.Range("I107").Value = 1200
$ws.Range("K107").Value = 1200
$ws.Range("M107").Value = 720
$ws.Range("H113").Value = 501703
$ws.Range("I113").Value = 3400
$ws.Range("J113").Value = 1000006
$ws.Range("K113").Value = 3400
$ws.Range("L113").Value = 1000006
$ws.Range("M113").Value = -146
$ws.Range("N113").Value = -1006514
$ws.Range("H122").Value = 12090.786
$ws.Range("I122").Value = 14597.637
$ws.Range("K122").Value = 43792.911
$ws.Range("M122").Value = -41342.911
$ws.Range("H137").Value = 3786.9546
$ws.Range("I137").Value = 2211.4443
$ws.Range("K137").Value = 6634.3329
$ws.Range("M137").Value = -4084.3329
$ws.Range("H138").Value = 3148.8298
$ws.Range("I138").Value = 2556.1428
$ws.Range("J138").Value = 3400.2727
$ws.Range("K138").Value = 7668.428400000001
$ws.Range("L138").Value = 10200.8181
$ws.Range("M138").Value = -2528.428400000001
$ws.Range("N138").Value = -20480.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5112756
$ws.Range("I2").Value = 9202162
$ws.Range("K2").Value = 9202162
$ws.Range("M2").Value = -9202049
$ws.Range("H61").Value = 17550536
$ws.Range("I61").Value = 25644932
$ws.Range("J61").Value = 12677.833
$ws.Range("K61").Value = 25644932
$ws.Range("L61").Value = 12677.833
$ws.Range("M61").Value = -25644720
$ws.Range("N61").Value = -13101.833
$ws.Range("H116").Value = 5112756
$ws.Range("I116").Value = 9202162
$ws.Range("K116").Value = 9202162
$ws.Range("M116").Value = -9199868
$ws.Range("H122").Value = 1256.4166
$ws.Range("I122").Value = 884.125
$ws.Range("J122").Value = 2001
$ws.Range("K122").Value = 2652.375
$ws.Range("L122").Value = 6003
$ws.Range("M122").Value = -202.375
$ws.Range("N122").Value = -10903
$ws.Range("H132").Value = 6903472.5
$ws.Range("I132").Value = 12504728
$ws.Range("J132").Value = 9620.538
$ws.Range("K132").Value = 37514184
$ws.Range("L132").Value = 28861.614
$ws.Range("M132").Value = -37511654
$ws.Range("N132").Value = -33921.614
$ws.Range("H133").Value = 119999.5
$ws.Range("J133").Value = 119999.5
$ws.Range("L133").Value = 119999.5
$ws.Range("N133").Value = -125059.5
$ws.Range("H136").Value = 17550536
$ws.Range("I136").Value = 25644932
$ws.Range("J136").Value = 12677.833
$ws.Range("K136").Value = 76934796
$ws.Range("L136").Value = 38033.499
$ws.Range("M136").Value = -76932246
$ws.Range("N136").Value = -43133.499
$ws.Range("H140").Value = 80428
$ws.Range("J140").Value = 80428
$ws.Range("L140").Value = 80428
$ws.Range("N140").Value = -90788

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5112756
$ws.Range("I3").Value = 9202162
$ws.Range("K3").Value = 9202162
$ws.Range("M3").Value = -9202048
$ws.Range("H107").Value = 1151.4445
$ws.Range("I107").Value = 1117.0834
$ws.Range("J107").Value = 1426.3334
$ws.Range("K107").Value = 1117.0834
$ws.Range("L107").Value = 1426.3334
$ws.Range("M107").Value = 802.9166
$ws.Range("N107").Value = -5266.3334
$ws.Range("H134").Value = 5032.8335
$ws.Range("I134").Value = 3353.6428
$ws.Range("J134").Value = 8391.214
$ws.Range("K134").Value = 10060.9284
$ws.Range("L134").Value = 25173.642
$ws.Range("M134").Value = -7525.928400000001
$ws.Range("N134").Value = -30243.642
$ws.Range("H135").Value = 104100.875
$ws.Range("J135").Value = 104100.875
$ws.Range("L135").Value = 104100.875
$ws.Range("N135").Value = -114240.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2270.4
$ws.Range("I16").Value = 1973
$ws.Range("J16").Value = 3162.6
$ws.Range("K16").Value = 1973
$ws.Range("L16").Value = 3162.6
$ws.Range("M16").Value = -1686
$ws.Range("N16").Value = -3736.6
$ws.Range("H31").Value = 5286.893
$ws.Range("J31").Value = 6852
$ws.Range("L31").Value = 6852
$ws.Range("N31").Value = -7442
$ws.Range("H34").Value = 5286.893
$ws.Range("J34").Value = 6852
$ws.Range("L34").Value = 6852
$ws.Range("N34").Value = -7256
$ws.Range("H99").Value = 6666.3335
$ws.Range("J99").Value = 7500
$ws.Range("L99").Value = 7500
$ws.Range("N99").Value = -10496
$ws.Range("H105").Value = 1321.3334
$ws.Range("I105").Value = 1185.6
$ws.Range("K105").Value = 1185.6
$ws.Range("M105").Value = 561.4000000000001
$ws.Range("H113").Value = 2270.4
$ws.Range("I113").Value = 1973
$ws.Range("J113").Value = 3162.6
$ws.Range("K113").Value = 1973
$ws.Range("L113").Value = 3162.6
$ws.Range("M113").Value = 197
$ws.Range("N113").Value = -7502.6
$ws.Range("H122").Value = 2246.8235
$ws.Range("I122").Value = 607.46155
$ws.Range("J122").Value = 7574.75
$ws.Range("K122").Value = 1822.38465
$ws.Range("L122").Value = 22724.25
$ws.Range("M122").Value = 627.61535
$ws.Range("N122").Value = -27624.25
$ws.Range("H126").Value = 6666.3335
$ws.Range("J126").Value = 7500
$ws.Range("L126").Value = 22500
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4831.4
$ws.Range("J63").Value = 4788
$ws.Range("L63").Value = 14364
$ws.Range("N63").Value = -15862
$ws.Range("H66").Value = 4831.4
$ws.Range("J66").Value = 4788
$ws.Range("L66").Value = 43092
$ws.Range("N66").Value = -50580
$ws.Range("H80").Value = 4098.5
$ws.Range("J80").Value = 4098.5
$ws.Range("L80").Value = 12295.5
$ws.Range("N80").Value = -14167.5
$ws.Range("H83").Value = 4098.5
$ws.Range("J83").Value = 4098.5
$ws.Range("L83").Value = 36886.5
$ws.Range("N83").Value = -46246.5
$ws.Range("H104").Value = 3991.25
$ws.Range("J104").Value = 3988.3333
$ws.Range("L104").Value = 11964.9999
$ws.Range("N104").Value = -17206.9999
$ws.Range("H126").Value = 1900
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
$ws.Range("H131").Value = 19616892
$ws.Range("I131").Value = 83334390
$ws.Range("J131").Value = 11507.538
$ws.Range("K131").Value = 250003170
$ws.Range("L131").Value = 34522.614
$ws.Range("M131").Value = -249998130
$ws.Range("N131").Value = -44602.614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H126").Value = 3380.6
$ws.Range("I126").Value = 3022.25
$ws.Range("K126").Value = 9066.75
$ws.Range("M126").Value = -6596.75
$ws.Range("H132").Value = 5053
$ws.Range("I132").Value = 3444.5557
$ws.Range("J132").Value = 6431.6665
$ws.Range("K132").Value = 10333.6671
$ws.Range("L132").Value = 19294.9995
$ws.Range("M132").Value = -7803.667099999999
$ws.Range("N132").Value = -24354.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 21000
$ws.Range("J101").Value = 21000
$ws.Range("L101").Value = 21000
$ws.Range("N101").Value = -27490
$ws.Range("H132").Value = 4008.2954
$ws.Range("I132").Value = 2625.8276
$ws.Range("J132").Value = 6681.067
$ws.Range("K132").Value = 7877.4828
$ws.Range("L132").Value = 20043.201
$ws.Range("M132").Value = -5347.4828
$ws.Range("N132").Value = -25103.201
$ws.Range("H136").Value = 1240335.5
$ws.Range("J136").Value = 8027.5557
$ws.Range("L136").Value = 24082.6671
$ws.Range("N136").Value = -29182.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 25000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 25000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 25000
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = -26248
$ws.Range("H66").Value = 25000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 25000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 75000
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = -81240
$ws.Range("H101").Value = 26996.857
$ws.Range("J101").Value = 26996.857
$ws.Range("L101").Value = 26996.857
$ws.Range("N101").Value = -33486.857
$ws.Range("H103").Value = 53797.2
$ws.Range("J103").Value = 53797.2
$ws.Range("L103").Value = 53797.2
$ws.Range("N103").Value = -56141.2
$ws.Range("H132").Value = 4702.3423
$ws.Range("I132").Value = 3536.5806
$ws.Range("J132").Value = 9865
$ws.Range("K132").Value = 10609.7418
$ws.Range("L132").Value = 29595
$ws.Range("M132").Value = -8079.7418
$ws.Range("N132").Value = -34655
$ws.Range("H136").Value = 2559.1714
$ws.Range("I136").Value = 1371.68
$ws.Range("K136").Value = 4115.04
$ws.Range("M136").Value = -1565.04
